{"js": "// The document carries a hidden \"_GoBack\" bookmark (Word's \"last edit\n// position\" marker) sitting right after the \"Excel:\" run. As further\n// edits were made at the end of the document, Word moved this bookmark\n// to the new last-edit location: the (empty) final paragraph. Replicate\n// that by deleting the bookmark from its old spot and re-inserting it,\n// collapsed, at the end of the last paragraph.\n\nconst doc = context.document;\n\n// Remove the existing \"_GoBack\" bookmark (currently right after \"Excel:\").\ndoc.deleteBookmark(\"_GoBack\");\n\nconst paragraphs = doc.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst endOfDoc = lastParagraph.getRange(\"End\");\nendOfDoc.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document carries a hidden \"_GoBack\" bookmark (tracks the last edit\n# position) sitting right after the \"Excel:\" run. Word re-stamps this\n# bookmark at the new last-edit location whenever the document is edited\n# further down. Here it needs to move to the very end of the document,\n# i.e. into the (empty) final paragraph, leaving the \"Excel:\" paragraph\n# bookmark-free.\n$bm = $d.Bookmarks(\"_GoBack\")\n$bm.Delete()\n\n$lastPara = $d.Paragraphs($d.Paragraphs.Count)\n$r = $lastPara.Range\n$r.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $r)\n"}
